$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D7").Value = "Objective`nid=""Par1"">To investigate the occurrence of olfactory and gustatory dysfunctions in patients with laboratory-confirmed COVID-19 infection.`n`n`nMethods`nid=""Par2"">Patients with laboratory-confirmed COVID-19 infection were recruited from 12 European hospitals.`n`n The following epidemiological and clinical outcomes have been studied: age, sex, ethnicity, comorbidities, and general and otolaryngological symptoms.`n`n Patients completed olfactory and gustatory questionnaires based on the smell and taste component of the National Health and Nutrition Examination Survey, and the short version of the Questionnaire of Olfactory Disorders-Negative Statements (sQOD-NS).`n`n`nResults`nid=""Par3"">A total of 417 mild-to-moderate COVID-19 patients completed the study (263 females).`n`n The most prevalent general symptoms consisted of cough, myalgia, and loss of appetite.`n`n Face pain and nasal obstruction were the most disease-related otolaryngological symptoms.`n`n 85.6% and 88.0% of patients reported olfactory and gustatory dysfunctions, respectively.`n`n There was a significant association between both disorders (p &lt; 0.001).`n`n Olfactory dysfunction (OD) appeared before the other symptoms in 11.8% of cases.`n`n The sQO-NS scores were significantly lower in patients with anosmia compared with normosmic or hyposmic individuals (p = 0.001).`n`n Among the 18.2% of patients without nasal obstruction or rhinorrhea, 79.7% were hyposmic or anosmic.`n`n The early olfactory recovery rate was 44.0%.`n`n Females were significantly more affected by olfactory and gustatory dysfunctions than males (p = 0.001).`n`n`nConclusion`nid=""Par4"">Olfactory and gustatory disorders are prevalent symptoms in European COVID-19 patients, who may not have nasal symptoms.`n`n The sudden anosmia or ageusia need to be recognized by the international scientific community as important symptoms of the COVID-19 infection.`n`n`nElectronic supplementary material`nThe online version of this article (10.1007/s00405-020-05965-1) contains supplementary material, which is available to authorized users.`n`n`n"
$ws.Range("E7").Value = "[Jerome R.%Lechien%Jerome.Lechien@umons.ac.be%0,  Carlos M.%Chiesa-Estomba%NULL%0,  Daniele R.%De Siati%NULL%0,  Mihaela%Horoi%NULL%0,  Serge D.%Le Bon%NULL%0,  Alexandra%Rodriguez%NULL%0,  Didier%Dequanter%NULL%0,  Serge%Blecic%NULL%0,  Fahd%El Afia%NULL%0,  Lea%Distinguin%NULL%0,  Younes%Chekkoury-Idrissi%NULL%0,  Stéphane%Hans%NULL%0,  Irene Lopez%Delgado%NULL%0,  Christian%Calvo-Henriquez%NULL%0,  Philippe%Lavigne%NULL%0,  Chiara%Falanga%NULL%0,  Maria Rosaria%Barillari%NULL%0,  Giovanni%Cammaroto%NULL%0,  Mohamad%Khalife%NULL%0,  Pierre%Leich%NULL%0,  Christel%Souchay%NULL%0,  Camelia%Rossi%NULL%0,  Fabrice%Journe%NULL%0,  Julien%Hsieh%NULL%0,  Myriam%Edjlali%NULL%0,  Robert%Carlier%NULL%0,  Laurence%Ris%NULL%0,  Andrea%Lovato%NULL%0,  Cosimo%De Filippis%NULL%0,  Frederique%Coppee%NULL%0,  Nicolas%Fakhry%NULL%0,  Tareck%Ayad%NULL%0,  Sven%Saussez%NULL%0]"
